# edit.ps1
# Applies country/provincia data updates to the "Pais" worksheet as described
# in the commit "Update countries & provincias Spain".
#
# Summary of changes:
#  - Updates the "last updated" timestamp string (A1)
#  - Refreshes case counts for several countries (rows resorted by total cases
#    causes Pakistan/Chile and Santa Sede/Surinam/Gambia to swap positions)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 28 de Abril de 2020 a las 19:22'
$ws.Range("B4").Value = 1019823
$ws.Range("C4").Value = 9467
$ws.Range("D4").Value = 139927
$ws.Range("E4").Value = 822295
$ws.Range("F4").Value = 14145
$ws.Range("G4").Value = 804
$ws.Range("H4").Value = 57601
$ws.Range("B8").Value = 159137
$ws.Range("C8").Value = 379
$ws.Range("E8").Value = 35563
$ws.Range("G8").Value = 48
$ws.Range("H8").Value = 6174
$ws.Range("D10").Value = 38809
$ws.Range("E10").Value = 72852
$ws.Range("A30").Value = 'Pakistan'
$ws.Range("B30").Value = 14514
$ws.Range("C30").Value = 599
$ws.Range("D30").Value = 3233
$ws.Range("E30").Value = 10969
$ws.Range("F30").Value = 111
$ws.Range("G30").Value = 20
$ws.Range("H30").Value = 312
$ws.Range("A31").Value = 'Chile'
$ws.Range("B31").Value = 14365
$ws.Range("C31").Value = 552
$ws.Range("D31").Value = 7710
$ws.Range("E31").Value = 6448
$ws.Range("F31").Value = 426
$ws.Range("H31").Value = 207
$ws.Range("B44").Value = 7619
$ws.Range("C44").Value = 20
$ws.Range("E44").Value = 7381
$ws.Range("F44").Value = 44
$ws.Range("B65").Value = 2566
$ws.Range("C65").Value = 32
$ws.Range("E65").Value = 1851
$ws.Range("F65").Value = 40
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = 138
$ws.Range("D68").Value = 992
$ws.Range("E68").Value = 939
$ws.Range("B105").Value = 619
$ws.Range("C105").Value = 31
$ws.Range("E105").Value = 478
$ws.Range("A202").Value = 'Santa Sede'
$ws.Range("C202").Value = 1
$ws.Range("D202").Value = 2
$ws.Range("E202").Value = 8
$ws.Range("H202").Value = 0
$ws.Range("A203").Value = 'Surinam'
$ws.Range("D203").Value = 7
$ws.Range("E203").Value = 2
$ws.Range("A204").Value = 'Gambia'
$ws.Range("B204").Value = 10
$ws.Range("D204").Value = 8
$ws.Range("E204").Value = 1
$ws.Range("H204").Value = 1
